$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-23: Name, Code, Width, Height, Thickness, Color, Qty, Grain
$rows = @(
    @("柜体侧板（L）", "CB(L)-HS00-2434-574-16", 2434, 574, 16, "HS00", 1, "mixed"),
    @("柜体侧板（R）", "CB(R)-HS00-2434-574-16", 2434, 574, 16, "HS00", 1, "mixed"),
    @("顶底板（双门柜体）", "DD-HS00-704-554-16", 704, 554, 16, "HS00", 2, "mixed"),
    @("顶底板（双门柜体）", "DD-HS00-896-554-16", 896, 554, 16, "HS00", 1, "mixed"),
    @("顶底板（双门柜体）", "DD-HS00-832-554-16", 832, 554, 16, "HS00", 2, "mixed"),
    @("顶底板（双门柜体）", "DD-HS00-736-554-16", 736, 554, 16, "HS00", 2, "mixed"),
    @("层隔板（双门柜体）", "CG-HS00-896-554-25", 896, 554, 25, "HS00", 3, "mixed"),
    @("层隔板（双门柜体）", "CG-HS00-704-554-25", 704, 554, 25, "HS00", 3, "mixed"),
    @("层隔板（单门柜体）", "CG-HS00-336-554-25", 336, 554, 25, "HS00", 1, "mixed"),
    @("后背板（双门柜体）为(1+1)组合", "HB-HS00-2320-464-12", 2320, 464, 12, "HS00", 1, "mixed"),
    @("后背板（双门柜体）为(1+1)组合", "HB-HS00-2256-336-12", 2256, 336, 12, "HS00", 1, "mixed"),
    @("门板（L/R）", "MB(R)-(门板花色)-2320-397-16", 2320, 397, 16, "(门板花色)", 2, "fixed"),
    @("门板（L/R）", "MB(L)-(门板花色)-2288-397-16", 2288, 397, 16, "(门板花色)", 2, "fixed"),
    @("门板（L/R）", "MB(L)-(门板花色)-2320-397-16", 2320, 397, 16, "(门板花色)", 2, "fixed"),
    @("双抽屉组件", "CTF2抽屉吊板(R)-HS00-534-371-16", 534, 371, 16, "HS00", 3, "mixed"),
    @("双抽屉组件", "抽屉拉板-HS00-438-106-12", 438, 106, 12, "HS00", 2, "mixed"),
    @("双抽屉组件", "800抽屉面板-HS00-800-170-16", 800, 170, 16, "HS00", 3, "mixed"),
    @("单抽屉组件", "526抽屉后板-HS00-526-138-12", 526, 138, 12, "HS00", 1, "mixed"),
    @("收口条", "TSB50-HS98-2434-50-16", 2434, 50, 16, "HS98", 2, "mixed"),
    @("底支撑(双门柜体）", "DC-HS03-672-82-16", 672, 82, 16, "HS03", 1, "fixed"),
    @("底支撑(双门柜体）", "DC-HS97-832-82-16", 832, 82, 16, "HS97", 2, "mixed"),
    @("底支撑(单门柜体）", "DC-HS02-368-82-16", 368, 82, 16, "HS02", 1, "fixed")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Remove the two trailing rows (24 and 25) that no longer exist in the updated table
$ws.Rows(24).Delete()
$ws.Rows(24).Delete()
